# Update cryptocurrency price/volume snapshot values (Price column D, Volume(1h) column E).
# Values that look like plain numbers are written with a leading apostrophe so Excel
# keeps them as literal text (matching the original inlineStr cells) instead of coercing
# them to numeric cells; the style is then reset to "Normal" so no stray number format
# is left attached to the cell.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "49.888.46"
$ws.Range("E2").Value = "  +3.13%  "
$ws.Range("D3").Value = "2.622.14"
$ws.Range("E3").Value = "  +4.40%  "
$ws.Range("E4").Value = "  -0.11%  "
$ws.Range("D5").Value = "'326.96"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +1.67%  "
$ws.Range("D6").Value = "'109.91"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +1.17%  "
$ws.Range("E7").Value = "  +0.93%  "
$ws.Range("E8").Value = "  -0.10%  "
$ws.Range("E9").Value = "  +3.30%  "
$ws.Range("D10").Value = "'40.32"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +0.75%  "
$ws.Range("D11").Value = "'20.62"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +1.40%  "
$ws.Range("E13").Value = "  +0.55%  "
$ws.Range("E14").Value = "  +1.29%  "
$ws.Range("D15").Value = "3.033.19"
$ws.Range("E15").Value = "  +4.49%  "
$ws.Range("D16").Value = "2.606.11"
$ws.Range("E16").Value = "  +3.80%  "
$ws.Range("D17").Value = "'0.873"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  +3.32%  "
$ws.Range("D18").Value = "49.776.83"
$ws.Range("E18").Value = "  +3.24%  "
$ws.Range("D19").Value = "'3.06"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +11.08%  "
$ws.Range("D20").Value = "'13.35"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +1.58%  "
$ws.Range("E21").Value = "  +0.28%  "
$ws.Range("D22").Value = "0.0₃0958"
$ws.Range("E22").Value = "  +1.43%  "
$ws.Range("D23").Value = "'72.65"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +0.50%  "
$ws.Range("D24").Value = "'278.26"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -0.22%  "
$ws.Range("E25").Value = "  +1.48%  "
$ws.Range("D26").Value = "'26.51"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +2.78%  "
$ws.Range("E27").Value = "  -0.01%  "
$ws.Range("D28").Value = "'2.22"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +0.64%  "
$ws.Range("E29").Value = "  +1.56%  "
$ws.Range("D30").Value = "'36.59"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +3.49%  "
$ws.Range("D31").Value = "'0.143"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +1.91%  "
$ws.Range("D32").Value = "'49.78"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +0.37%  "
$ws.Range("E33").Value = "  -0.13%  "
$ws.Range("E34").Value = "  +1.62%  "
$ws.Range("E35").Value = "  -0.17%  "
$ws.Range("E36").Value = "  +0.85%  "
$ws.Range("E37").Value = "  +5.04%  "
$ws.Range("D38").Value = "'4.75"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +1.50%  "
$ws.Range("D39").Value = "'3.11"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +6.26%  "
$ws.Range("D40").Value = "'124.01"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +1.36%  "
$ws.Range("E41").Value = "  +0.56%  "
$ws.Range("D42").Value = "'22.55"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +4.31%  "
$ws.Range("D43").Value = "'2.23"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +0.28%  "
$ws.Range("E44").Value = "  +3.88%  "
$ws.Range("D45").Value = "'3.36"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +4.58%  "
$ws.Range("D46").Value = "2.049.87"
$ws.Range("D47").Value = "'2.35"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +17.60%  "
$ws.Range("E48").Value = "  +8.41%  "
$ws.Range("D49").Value = "'9.03"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -0.21%  "
$ws.Range("E50").Value = "  +3.12%  "
$ws.Range("D51").Value = "'81.48"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +0.89%  "
